$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (shifts old D:K to F:M)
$ws.Columns("D:E").Insert()

# Copy number formats (styles) from the now-shifted F:G columns into the
# new D:E columns so the inserted cells carry the correct style (dates in
# row 7/38/80, "#,##0" for data rows, etc.) instead of Excel's default.
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 271400
$ws.Range("E8").Value = 263200
$ws.Range("D9:E9").Value = "NA"
$ws.Range("D10:E10").Value = "NA"
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = -5500
$ws.Range("E15").Value = -5600
$ws.Range("D17").Value = 62000
$ws.Range("E17").Value = 55900
$ws.Range("D18").Value = 209400
$ws.Range("E18").Value = 207300
$ws.Range("D20").Value = -104900
$ws.Range("E20").Value = -105700
$ws.Range("D21").Value = 116800
$ws.Range("E21").Value = 114000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 104500
$ws.Range("E23").Value = 101700
$ws.Range("D24").Value = 8300
$ws.Range("E24").Value = 17800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 96200
$ws.Range("E26").Value = 83900
$ws.Range("D27").Value = 94500
$ws.Range("E27").Value = 82300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29:E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 104900
$ws.Range("E32").Value = 105700
$ws.Range("D33").Value = 94500
$ws.Range("E33").Value = 82300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 94500
$ws.Range("E35").Value = 82300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 383400
$ws.Range("E41").Value = 339600
$ws.Range("D42").Value = 111100
$ws.Range("E42").Value = 108100
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2695400
$ws.Range("E48").Value = 1867300
$ws.Range("D49").Value = 887100
$ws.Range("E49").Value = 892600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 23000
$ws.Range("E52").Value = 59600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 28235900
$ws.Range("E54").Value = 28098200
$ws.Range("D57").Value = 12300
$ws.Range("E57").Value = 16000
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 225000
$ws.Range("E61").Value = 215900
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 25154600
$ws.Range("E66").Value = 25119300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1243600
$ws.Range("E72").Value = 1170900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3081300
$ws.Range("E76").Value = 2978900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 94500
$ws.Range("E81").Value = 82300
$ws.Range("D83").Value = 12300
$ws.Range("E83").Value = 12400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 119100
$ws.Range("E89").Value = 109900
$ws.Range("D91").Value = -18100
$ws.Range("E91").Value = -14500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -101600
$ws.Range("E94").Value = 20100
$ws.Range("D96").Value = -23600
$ws.Range("E96").Value = -23500
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 26200
$ws.Range("E100").Value = -145500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 43800
$ws.Range("E102").Value = -15500
